$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as text (matching the existing rows above),
# so force a text format before assigning the numeric-looking string,
# otherwise Excel would auto-detect it as a number.
$ws.Cells.Item(122, 1).NumberFormat = "@"
$ws.Cells.Item(122, 1).Value = "20210414"
$ws.Cells.Item(122, 2).Value = 658.0

$ws.Cells.Item(123, 1).NumberFormat = "@"
$ws.Cells.Item(123, 1).Value = "20210428"
$ws.Cells.Item(123, 2).Value = 750.0
